$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "updated age split model": fill in the previously-missing DIC pw/non-pw
# split for External causes, sex=1 (male), age=45 (row 27).
$ws.Range("D27").Value = 111180.4
$ws.Range("E27").Value = 111186.12
$ws.Range("F27").Formula = "=IF(D27<E27,""non-pw"",""pw"")"

# Reflect the updated scroll position / active selection on the sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G15").Select()
